$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-82 down to 13-83
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new record's data
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 44462
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 100112040
$ws.Cells.Item(12, 7).Value = "Cilantro"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 3200
$ws.Cells.Item(12, 11).Value = 1500
$ws.Cells.Item(12, 12).Value = 2000
$ws.Cells.Item(12, 13).Value = 1750
$ws.Cells.Item(12, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 16).Value = 1167
$ws.Cells.Item(12, 17).Value = 1.5
$ws.Cells.Item(12, 18).Value = "Hortaliza"
